$d = $word.ActiveDocument
$last = $d.Paragraphs.Last
$r = $last.Range
Write-Host "last para start=$($r.Start) end=$($r.End) text=[$($r.Text)]"

$insPoint = $d.Range($r.End, $r.End)
Write-Host "insPoint start=$($insPoint.Start) end=$($insPoint.End)"

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
       '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
       '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>' +
       '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
       '<w:body>' +
       '<w:p><w:pPr><w:pStyle w:val="a3"/><w:widowControl/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="10"/></w:numPr><w:tabs><w:tab w:val="left" w:pos="709"/><w:tab w:val="left" w:pos="851"/><w:tab w:val="left" w:pos="2410"/></w:tabs><w:suppressAutoHyphens w:val="0"/><w:autoSpaceDE/><w:spacing w:before="40" w:after="60" w:line="360" w:lineRule="auto"/><w:ind w:left="1701" w:right="567" w:firstLine="567"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="44"/><w:szCs w:val="44"/></w:rPr></w:pPr>' +
       '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:color w:val="000000"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t>незарегистрированный клиент не может записаться на занятие</w:t></w:r>' +
       '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:color w:val="000000"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t>.</w:t></w:r>' +
       '</w:p>' +
       '</w:body></w:document>' +
       '</pkg:xmlData></pkg:part></pkg:package>'

try {
  $insPoint.InsertXML($xml)
  Write-Host "InsertXML ok"
} catch {
  Write-Host "InsertXML failed: $_"
}
Write-Host "Paragraphs count: $($d.Paragraphs.Count)"
Write-Host "Last para text: [$($d.Paragraphs.Last.Range.Text)]"
